# UI Validations added V.12
# Append two new regression test rows (NC_OP_22, NC_OP_23) to the TestCases
# sheet, matching the formatting already used by the rest of the table, then
# extend the Yes/No list validation to cover the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# --- Row 26 -----------------------------------------------------------
# Clone row 25's formatting (values + formats) down into row 26, then
# overwrite with the new row's content.
$ws.Range("A25:E25").Copy()
$ws.Range("A26:E26").PasteSpecial(-4122)
$ws.Range("A25:E25").Copy()
$ws.Paste($ws.Range("A26:E26"))

$ws.Cells.Item(26, 1).Value = "SAP Regression Automation"
$ws.Cells.Item(26, 2).Value = "Yes"
$ws.Cells.Item(26, 3).Value = "NC_OP_22"
$ws.Cells.Item(26, 4).Value = 5397361
$ws.Cells.Item(26, 5).Value = "Create PP Society discount with Rejected"

# --- Row 27 -----------------------------------------------------------
# Clone row 26's formatting down into row 27, then overwrite with the new
# row's content.
$ws.Range("A26:E26").Copy()
$ws.Range("A27:E27").PasteSpecial(-4122)
$ws.Range("A26:E26").Copy()
$ws.Paste($ws.Range("A27:E27"))

$ws.Cells.Item(27, 1).Value = "SAP Regression Automation"
$ws.Cells.Item(27, 2).Value = "Yes"
$ws.Cells.Item(27, 3).Value = "NC_OP_23"
$ws.Cells.Item(27, 4).Value = 5397362
$ws.Cells.Item(27, 5).Value = "Create PP Society discount with Withdrawn"

$excel.CutCopyMode = 0

# --- Data validation ----------------------------------------------------
# Extend the Yes/No dropdown list down through the two new rows.
$ws.Range("B2:B27").Validation.Delete()
$ws.Range("B2:B27").Validation.Add(3, 1, 1, '"Yes,No"')

# --- View state -----------------------------------------------------------
# Mirror the author's final scroll position / selection.
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("B31").Select()
